$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '30.283.34'
$ws.Range('E2').Value2 = '  +0.10%  '

$ws.Range('D3').Value2 = '1.868.10'
$ws.Range('E3').Value2 = '  +0.24%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value2 = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '235.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  -0.86%  '

$ws.Range('E6').Value2 = '  +0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.4690'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.2861'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -0.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.06585'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +0.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '21.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -1.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.07957'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  +0.64%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '96.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  -1.39%  '

$ws.Range('D13').Value2 = '1.877.17'
$ws.Range('E13').Value2 = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '0.6893'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +0.39%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '5.103'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  -1.82%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '268.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  -3.60%  '

$ws.Range('D17').Value2 = '30.354.76'
$ws.Range('E17').Value2 = '  +0.31%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '14.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  +3.55%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.000007764'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +5.39%  '

$ws.Range('E20').Value2 = '  +0.03%  '

$ws.Range('D21').Value2 = '2.121.95'
$ws.Range('E21').Value2 = '  +0.31%  '

$ws.Range('E22').Value2 = '  -0.02%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '5.245'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  -1.90%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '6.208'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  +0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '9.373'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +1.18%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '167.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value2 = '  -0.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '18.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -1.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '1.946'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  -0.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '1.363'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  -1.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '0.09863'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +0.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '4.335'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  -1.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '1.457'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -1.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '4.048'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  -0.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.04709'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  -0.63%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '1.134'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  -0.47%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.7014'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  -0.49%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '2.733'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  +0.88%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.01874'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  -0.36%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '2.793'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  +6.36%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '6.247'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  -0.81%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '72.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  -4.79%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '1.955'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  -0.21%  '

$ws.Range('B43').Value2 = 'TrustWalletToken'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.8413'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  -1.30%  '

$ws.Range('B44').Value2 = 'TheSandbox'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.4168'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -0.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +0.00%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '102.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  -0.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '7.073'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  -2.07%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '9.092'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -2.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '917.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  -3.82%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '34.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +0.68%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '0.05691'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  +0.78%  '
